$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.539.37"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").Value = "2.642.51"
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.13"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("E8").Value = "  +1.23%  "

$ws.Range("D9").Value = "2.666.89"
$ws.Range("E9").Value = "  +1.70%  "

$ws.Range("E10").Value = "  +2.94%  "

$ws.Range("E11").Value = "  +1.89%  "

$ws.Range("E12").Value = "  -0.30%  "

$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("D14").Value = "3.112.25"
$ws.Range("E14").Value = "  +1.29%  "

$ws.Range("D15").Value = "59.514.72"
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.30"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").Value = "2.665.21"
$ws.Range("E18").Value = "  +1.75%  "

$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.54"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.59%  "

$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("E25").Value = "  +1.31%  "

$ws.Range("D26").Value = "2.764.85"
$ws.Range("E26").Value = "  +1.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.162"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("E29").Value = "  +1.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("E31").Value = "  -0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.53"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.11"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.69"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("E36").Value = "  +14.94%  "

$ws.Range("E37").Value = "  +2.86%  "

$ws.Range("E38").Value = "  +2.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.874"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.57"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.74"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.32%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "284.95"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.620"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0999"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.82"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0547"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("E49").Value = "  +1.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.76"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.18%  "

$ws.Range("E51").Value = "  -1.24%  "
